# Data-dictionary update: add "csv abbreviation" / "Sampling frequency"
# columns to the "empatica" sheet and a new "Interbeat Interval" (IBI) row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("empatica")
$ws.Activate()

# Insert two new columns before the old "Description" column (D), which
# pushes the existing "Description"/"Use" columns from D/E to F/G.
$ws.Columns("D:E").Insert()

# Match the width Excel gave the other narrow columns (B/C).
$ws.Columns("D:E").ColumnWidth = $ws.Columns("C").ColumnWidth

# --- Header row -----------------------------------------------------
$ws.Cells.Item(2, 4).Value = "csv abbreviation"
$ws.Cells.Item(2, 5).Value = "Sampling frequency"

# --- Existing rows: fill in the two new columns ----------------------
# Row 3: Electrodermal activity / EDA
$ws.Cells.Item(3, 4).Value = "EDA"
$ws.Cells.Item(3, 5).Value = "4 Hz"

# Row 4: Photoplethysmography sensors / PPG / BVP
$ws.Cells.Item(4, 4).Value = "BVP"
$ws.Cells.Item(4, 5).Value = "64 Hz"

# Row 5: Heart rate (beats per minute) / HR (BPM)
$ws.Cells.Item(5, 4).Value = "HR"
$ws.Cells.Item(5, 5).Value = "1Hz"

# Row 6: temperature signal
$ws.Cells.Item(6, 4).Value = "TEMP"
$ws.Cells.Item(6, 5).Value = "4 Hz"

# Row 7: 3D- acceleration meter
$ws.Cells.Item(7, 4).Value = "ACC"
$ws.Cells.Item(7, 5).Value = "32 Hz"

# --- New row 8: Interbeat Interval / IBI -----------------------------
$ws.Cells.Item(8, 1).Value = "Interbeat Interval"
$ws.Cells.Item(8, 2).Value = "IBI"
$ws.Cells.Item(8, 3).Value = "seconds"
$ws.Cells.Item(8, 4).Value = "IBI"
$ws.Cells.Item(8, 5).Value = "1 Hz"
$ws.Cells.Item(8, 7).Value = "DO NOT USE"

# Wrap text in the new columns, like the neighbouring description/use columns.
$ws.Range("D3:E8").WrapText = $true
$ws.Range($ws.Cells.Item(2,4), $ws.Cells.Item(2,5)).WrapText = $true

# Row heights: header wraps onto two lines now that it has more columns.
$ws.Rows(2).RowHeight = 34
$ws.Rows(8).RowHeight = 17

$ws.Range("G7").Select()

Write-Host "empatica sheet updated"
